# Update workbook "prod_AGV" — refreshed AGV production data via Shiny app
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "pro": refreshed production values (B2:B26)
# ---------------------------------------------------------------------------
$wsPro = $wb.Worksheets.Item("pro")
$arr = New-Object 'object[,]' 25,1
    $arr[0,0] = 1268339.4477639399
    $arr[1,0] = 1181742.0247480848
    $arr[2,0] = 1255521.6920803196
    $arr[3,0] = 1341359.1348832971
    $arr[4,0] = 1418177.9615245394
    $arr[5,0] = 1456786.3016296215
    $arr[6,0] = 1412661.8168779463
    $arr[7,0] = 1448417.6911906502
    $arr[8,0] = 1484911.6425653687
    $arr[9,0] = 1537487.4391951556
    $arr[10,0] = 1488682.1449741309
    $arr[11,0] = 1501662.42919933
    $arr[12,0] = 1593615.797616252
    $arr[13,0] = 1625302.9323412769
    $arr[14,0] = 1504966.1298186895
    $arr[15,0] = 1449487.3038074195
    $arr[16,0] = 1420710.0361085436
    $arr[17,0] = 1491774.0410095686
    $arr[18,0] = 1839728.8443069814
    $arr[19,0] = 1891618
    $arr[20,0] = 1884039
    $arr[21,0] = 2035994.6198774669
    $arr[22,0] = 2151223.5797187076
    $arr[23,0] = 2213602.0023530093
    $arr[24,0] = 2258851.356632086
$wsPro.Range("B2:B26").Value = $arr

# ---------------------------------------------------------------------------
# 2) Sheet "ind": quarterly figures, previously formulas referencing pro!,
#    now replaced with their refreshed literal values (B2:B101)
# ---------------------------------------------------------------------------
$wsInd = $wb.Worksheets.Item("ind")
$arr = New-Object 'object[,]' 100,1
    $arr[0,0] = 296815.1286447852
    $arr[1,0] = 632956.36367052561
    $arr[2,0] = 626505.23194503936
    $arr[3,0] = 679955.22614140587
    $arr[4,0] = 276549.71365822741
    $arr[5,0] = 589740.49581118638
    $arr[6,0] = 583729.82297385321
    $arr[7,0] = 633530.4536139745
    $arr[8,0] = 293815.53432571044
    $arr[9,0] = 626559.74796782504
    $arr[10,0] = 620173.81095853727
    $arr[11,0] = 673083.64300185547
    $arr[12,0] = 316340.56361062248
    $arr[13,0] = 674594.22886792955
    $arr[14,0] = 667718.72135830752
    $arr[15,0] = 724684.82469092519
    $arr[16,0] = 321235.7434833172
    $arr[17,0] = 685033.16864125268
    $arr[18,0] = 678051.26679006545
    $arr[19,0] = 735898.88629397925
    $arr[20,0] = 326307.65537126985
    $arr[21,0] = 695848.98830689269
    $arr[22,0] = 688756.85093018273
    $arr[23,0] = 747517.81222436635
    $arr[24,0] = 316424.14867733797
    $arr[25,0] = 674772.47348816507
    $arr[26,0] = 667895.14929797756
    $arr[27,0] = 724876.30449588026
    $arr[28,0] = 333162.00182239292
    $arr[29,0] = 710465.83827963448
    $arr[30,0] = 703224.72503349977
    $arr[31,0] = 763220.00608660234
    $arr[32,0] = 334065.98837462143
    $arr[33,0] = 712393.58382117108
    $arr[34,0] = 705132.82286922983
    $arr[35,0] = 765290.89237651567
    $arr[36,0] = 334989.62702592794
    $arr[37,0] = 714363.2373383143
    $arr[38,0] = 707082.40155180008
    $arr[39,0] = 767406.79843187612
    $arr[40,0] = 324355.87035826361
    $arr[41,0] = 691686.81924853066
    $arr[42,0] = 684637.10296496458
    $arr[43,0] = 743046.58993203379
    $arr[44,0] = 366939.75338007737
    $arr[45,0] = 782496.67746406433
    $arr[46,0] = 774521.42129979632
    $arr[47,0] = 840599.34589255904
    $arr[48,0] = 395516.07430650853
    $arr[49,0] = 843435.49909105152
    $arr[50,0] = 834839.15056074434
    $arr[51,0] = 906063.05337451329
    $arr[52,0] = 403380.43606246635
    $arr[53,0] = 860206.20024219714
    $arr[54,0] = 851438.92365355243
    $arr[55,0] = 924079.0281688089
    $arr[56,0] = 381715.25164920161
    $arr[57,0] = 814005.33303208067
    $arr[58,0] = 805708.93863581272
    $arr[59,0] = 874447.61135263159
    $arr[60,0] = 379497.4362782886
    $arr[61,0] = 809275.8559367751
    $arr[62,0] = 801027.66467343329
    $arr[63,0] = 869366.95674127678
    $arr[64,0] = 371963.11756702355
    $arr[65,0] = 793208.96946777613
    $arr[66,0] = 785124.53293853451
    $arr[67,0] = 852107.05692912615
    $arr[68,0] = 455373.53176948422
    $arr[69,0] = 971081.1443360073
    $arr[70,0] = 961183.82322856318
    $arr[71,0] = 1043186.7613583996
    $arr[72,0] = 530507.73802252766
    $arr[73,0] = 1131304.3586793416
    $arr[74,0] = 1119774.0323277195
    $arr[75,0] = 1215307.0182907304
    $arr[76,0] = 530920.11066465639
    $arr[77,0] = 1132183.7407015141
    $arr[78,0] = 1120644.4516320999
    $arr[79,0] = 1216251.6970017296
    $arr[80,0] = 527067.04315789125
    $arr[81,0] = 1123967.0988841159
    $arr[82,0] = 1112511.5543534216
    $arr[83,0] = 1207424.9454818098
    $arr[84,0] = 561667.63844664965
    $arr[85,0] = 1197752.6470628874
    $arr[86,0] = 1185545.0755078054
    $arr[87,0] = 1286689.2865604304
    $arr[88,0] = 577434.69570488972
    $arr[89,0] = 1231375.7958340677
    $arr[90,0] = 1218825.5349970728
    $arr[91,0] = 1322809.0525324726
    $arr[92,0] = 600234.5949825563
    $arr[93,0] = 1279996.4352358978
    $arr[94,0] = 1266950.629733644
    $arr[95,0] = 1375039.9167075262
    $arr[96,0] = 612205.1135233345
    $arr[97,0] = 1305523.4894713617
    $arr[98,0] = 1292217.5106003126
    $arr[99,0] = 1402462.4294297993
$wsInd.Range("B2:B101").Value = $arr

# ---------------------------------------------------------------------------
# 3) Sheet "conso": refreshed consumption values (B2:B26)
#    (VA sheet keeps its formula "=pro!B.. - conso!B.." and recalculates
#    automatically from the updated pro/conso values above.)
# ---------------------------------------------------------------------------
$wsConso = $wb.Worksheets.Item("conso")
$arr = New-Object 'object[,]' 25,1
    $arr[0,0] = 426293.43411799415
    $arr[1,0] = 397188.01924083667
    $arr[2,0] = 421985.76389978762
    $arr[3,0] = 450835.45185527805
    $arr[4,0] = 476654.52183515689
    $arr[5,0] = 489630.78989305423
    $arr[6,0] = 474800.69398198894
    $arr[7,0] = 486818.38613961462
    $arr[8,0] = 499084.7162199025
    $arr[9,0] = 516756.04112786177
    $arr[10,0] = 500352.38557762018
    $arr[11,0] = 504715.4308009877
    $arr[12,0] = 535621.0643051985
    $arr[13,0] = 546270.6341262823
    $arr[14,0] = 505825.21404671884
    $arr[15,0] = 487178.66442417534
    $arr[16,0] = 477506.02853767184
    $arr[17,0] = 501391.02148516086
    $arr[18,0] = 618339.71379275236
    $arr[19,0] = 635779
    $arr[20,0] = 621538
    $arr[21,0] = 672865.25258615182
    $arr[22,0] = 796526.13637147797
    $arr[23,0] = 819622.77980839321
    $arr[24,0] = 836377.10217498336
$wsConso.Range("B2:B26").Value = $arr

# ---------------------------------------------------------------------------
# 4) Recalculate so the "VA" formulas (pro - conso) pick up fresh values
# ---------------------------------------------------------------------------
$excel.CalculateFull()

# ---------------------------------------------------------------------------
# 5) Restore/update each sheet's selection; activate "pro" last so it
#    becomes the active tab of the workbook.
# ---------------------------------------------------------------------------
$wsInd.Range("D16").Select()
$wb.Worksheets.Item("VA").Range("D16").Select()
$wsConso.Range("D16").Select()
$wsPro.Range("D16").Select()
